$wb = $excel.ActiveWorkbook
$wsEvap = $wb.Worksheets.Item("evap")
$wsInflow = $wb.Worksheets.Item("Inflow")

# Update evap (sheet4) column B values, rows 4-34
$wsEvap.Range("B4").Value = 986
$wsEvap.Range("B5").Value = 986
$wsEvap.Range("B6").Value = 985
$wsEvap.Range("B7").Value = 984
$wsEvap.Range("B8").Value = 983
$wsEvap.Range("B9").Value = 983
$wsEvap.Range("B10").Value = 983
$wsEvap.Range("B11").Value = 984
$wsEvap.Range("B12").Value = 984
$wsEvap.Range("B13").Value = 984
$wsEvap.Range("B14").Value = 985
$wsEvap.Range("B15").Value = 985
$wsEvap.Range("B16").Value = 985
$wsEvap.Range("B17").Value = 984
$wsEvap.Range("B18").Value = 984
$wsEvap.Range("B19").Value = 983
$wsEvap.Range("B20").Value = 983
$wsEvap.Range("B21").Value = 983
$wsEvap.Range("B22").Value = 982
$wsEvap.Range("B23").Value = 982
$wsEvap.Range("B24").Value = 981
$wsEvap.Range("B25").Value = 981
$wsEvap.Range("B26").Value = 980
$wsEvap.Range("B27").Value = 980
$wsEvap.Range("B28").Value = 979
$wsEvap.Range("B29").Value = 978
$wsEvap.Range("B30").Value = 978
$wsEvap.Range("B31").Value = 977
$wsEvap.Range("B32").Value = 976
$wsEvap.Range("B33").Value = 976
$wsEvap.Range("B34").Value = 975

# Update Inflow (sheet5) column B values, rows 4-34
$wsInflow.Range("B4").Value = 4614
$wsInflow.Range("B5").Value = 5346
$wsInflow.Range("B6").Value = 4425
$wsInflow.Range("B7").Value = 4396
$wsInflow.Range("B8").Value = 7281
$wsInflow.Range("B9").Value = 7230
$wsInflow.Range("B10").Value = 13082
$wsInflow.Range("B11").Value = 14844
$wsInflow.Range("B12").Value = 11161
$wsInflow.Range("B13").Value = 13579
$wsInflow.Range("B14").Value = 16356
$wsInflow.Range("B15").Value = 7879
$wsInflow.Range("B16").Value = 9725
$wsInflow.Range("B17").Value = 8960
$wsInflow.Range("B18").Value = 4362
$wsInflow.Range("B19").Value = 5885
$wsInflow.Range("B20").Value = 9247
$wsInflow.Range("B21").Value = 9919
$wsInflow.Range("B22").Value = 5928
$wsInflow.Range("B23").Value = 8213
$wsInflow.Range("B24").Value = 8163
$wsInflow.Range("B25").Value = 9233
$wsInflow.Range("B26").Value = 3347
$wsInflow.Range("B27").Value = 6373
$wsInflow.Range("B28").Value = 6926
$wsInflow.Range("B29").Value = 5578
$wsInflow.Range("B30").Value = 5377
$wsInflow.Range("B31").Value = 6117
$wsInflow.Range("B32").Value = 5425
$wsInflow.Range("B33").Value = 5949
$wsInflow.Range("B34").Value = 5572
